# Form the consolidated report: recompute the "Absent" column (H) as the
# complement of the "Real" column (E) for every data row (rows 3-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 21; $row++) {
    $realValue = $ws.Cells.Item($row, 5).Value2   # Column E = "Real"
    $ws.Cells.Item($row, 8).Value = 1 - $realValue # Column H = "Absent"
}
